# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns
# for rows 2-51 per the source diff. D-column writes force a "@" (text)
# number format before assignment so numeric-looking strings such as
# "584.09" or "1.00" stay text (matching the original inlineStr cells)
# instead of being auto-coerced to a number by Excel; the format is
# reset back to Normal afterwards so no visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.322.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.419.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.71%  "
$ws.Range("E7").Value = "  +5.02%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.417.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.013.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.320.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.410.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "367.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("E24").Value = "  +6.94%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.869"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.713.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.95%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "336.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0286"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.72%  "
